$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F1").Value = "11_03_2024"
$ws.Range("F2").Value = 1405
$ws.Range("F3").Value = 1448
$ws.Range("F4").Value = 1548
$ws.Range("F5").Value = 3080
$ws.Range("F6").Value = 231

$ws.Range("F7").Select()
